# Import multiple courses/labs at once (commit: "Add multiple courses and
# labs at a time through importing").
#
# The "data" sheet is a lab-capacity table with header row (# / Lab name /
# Capacity) followed by 20 pre-numbered rows (1..20). We now fill in the
# real lab names + capacities for the first 11 rows and clear the leftover
# auto-numbering (column A) for the remaining, still-empty rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New lab data to import: Name, Capacity (row 2 already has "1" in col A).
$labs = @(
    @("PM1",      20),
    @("PM2",      20),
    @("PM4",      30),
    @("A5-301A",  30),
    @("A5-301B",  30),
    @("A5-302A",  30),
    @("A5-302B",  30),
    @("A5-303",   60),
    @("A5-304",   36),
    @("A5-203",   60),
    @("A5-204",   60)
)

$startRow = 2
for ($i = 0; $i -lt $labs.Count; $i++) {
    $r = $startRow + $i
    $name = $labs[$i][0]
    $cap = $labs[$i][1]
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = $cap
}

# The rows after the imported data no longer carry the placeholder
# auto-numbers in column A (they stay blank, ready for more imports).
$lastDataRow = $startRow + $labs.Count - 1
for ($r = $lastDataRow + 1; $r -le 21; $r++) {
    $ws.Cells.Item($r, 1).ClearContents()
}

# Match the cursor position left behind in the saved file.
$ws.Range("F17").Select() | Out-Null
